$p = $ppt.ActivePresentation

# --- 1) Re-colour the deck's theme (Design > Colors) from the custom
#        "Red Violet" (Integral) palette to the built-in "Office" palette.
#        PowerPoint stores RGB as 0x00BBGGRR in the ThemeColor.RGB property,
#        so each hex colour below is given in that reversed-byte form.
$s = $p.Slides.Item(1)
$officeThemeColors = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $s.ThemeColorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}

# --- 2) Re-style the table on slide 5 with the built-in table style that
#        replaces the previous custom "Table_0" style.
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{BC3790DD-833E-497C-BE74-681DDF63416E}")
    }
}
